$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Discover 100 Fortunes, an online slot game with expanding reels and Chinese aesthetics. Play now for free and read our expert review.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document (the one that repeats the H1 text in bold).
# ---------------------------------------------------------------------
$oldTitleText = "Play 100 Fortunes Free Online Slot | See Our Review"
$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 1; $i--) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq $oldTitleText) {
        $cand.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3) Replace the italic "meta description" paragraph's text with the
#    new image-prompt text (keep the italic run formatting intact).
# ---------------------------------------------------------------------
$newPromptText = 'Prompt: Create a feature image for "100 Fortunes" in a cartoon style featuring a happy Maya warrior with glasses. The image should depict the warrior standing in front of the game grid with a big smile on his face, holding a golden yo-yo and surrounded by Chinese cultural elements such as red paper lanterns and Ming vases. The cartoon style should be colorful and eye-catching, with a clear focus on the Maya warrior''s joyful expression. The image should convey the excitement and fun of playing this slot game while also highlighting the unique blend of Chinese and Maya cultures.'

$n2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n2)
$fullRange = $lastPara.Range
$textOnlyRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$promptXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>' + $newPromptText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$textOnlyRange.InsertXML($promptXml)

Write-Output "done"
